$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 432.05264
$ws.Range("I53").Value = 344.91666
$ws.Range("J53").Value = 581.4286
$ws.Range("K53").Value = 344.91666
$ws.Range("L53").Value = 581.4286
$ws.Range("M53").Value = 292.08334
$ws.Range("N53").Value = -1855.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5278.636
$ws.Range("I132").Value = 2093.261
$ws.Range("J132").Value = 12605
$ws.Range("K132").Value = 6279.782999999999
$ws.Range("L132").Value = 37815
$ws.Range("M132").Value = -3749.782999999999
$ws.Range("N132").Value = -42875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3173.8484
$ws.Range("I138").Value = 1928.5385
$ws.Range("J138").Value = 3479.302
$ws.Range("K138").Value = 5785.6155
$ws.Range("L138").Value = 10437.906
$ws.Range("M138").Value = -645.6154999999999
$ws.Range("N138").Value = -20717.906

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4801.4683
$ws.Range("I32").Value = 2868.0435
$ws.Range("K32").Value = 2868.0435
$ws.Range("M32").Value = -2581.0435

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5534.081
$ws.Range("I61").Value = 2460.12
$ws.Range("J61").Value = 11938.167
$ws.Range("K61").Value = 2460.12
$ws.Range("L61").Value = 11938.167
$ws.Range("M61").Value = -2248.12
$ws.Range("N61").Value = -12362.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5534.081
$ws.Range("I136").Value = 2460.12
$ws.Range("J136").Value = 11938.167
$ws.Range("K136").Value = 7380.36
$ws.Range("L136").Value = 35814.501
$ws.Range("M136").Value = -4830.36
$ws.Range("N136").Value = -40914.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7915.857
$ws.Range("I20").Value = 4003
$ws.Range("J20").Value = 10850.5
$ws.Range("K20").Value = 4003
$ws.Range("L20").Value = 10850.5
$ws.Range("M20").Value = -3756
$ws.Range("N20").Value = -11344.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2633.2632
$ws.Range("I86").Value = 1254.5834
$ws.Range("J86").Value = 4996.7144
$ws.Range("K86").Value = 1254.5834
$ws.Range("L86").Value = 4996.7144
$ws.Range("M86").Value = -131.5834
$ws.Range("N86").Value = -7242.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2633.2632
$ws.Range("I89").Value = 1254.5834
$ws.Range("J89").Value = 4996.7144
$ws.Range("K89").Value = 6272.916999999999
$ws.Range("L89").Value = 24983.572
$ws.Range("M89").Value = -656.9169999999995
$ws.Range("N89").Value = -36215.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5038.3335
$ws.Range("I134").Value = 3778.3635
$ws.Range("J134").Value = 6424.3
$ws.Range("K134").Value = 11335.0905
$ws.Range("L134").Value = 19272.9
$ws.Range("M134").Value = -8800.0905
$ws.Range("N134").Value = -24342.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2378.476
$ws.Range("I58").Value = 2485.2307
$ws.Range("J58").Value = 2205
$ws.Range("K58").Value = 2485.2307
$ws.Range("L58").Value = 2205
$ws.Range("M58").Value = -2282.2307
$ws.Range("N58").Value = -2611

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 19499.5
$ws.Range("J109").Value = 19499.5
$ws.Range("L109").Value = 19499.5
$ws.Range("N109").Value = -21579.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2378.476
$ws.Range("I136").Value = 2485.2307
$ws.Range("J136").Value = 2205
$ws.Range("K136").Value = 7455.6921
$ws.Range("L136").Value = 6615
$ws.Range("M136").Value = -4905.6921
$ws.Range("N136").Value = -11715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 442.85
$ws.Range("I122").Value = 214.25
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1928.25
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = 521.75
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 15626997
$ws.Range("I129").Value = 17857734
$ws.Range("J129").Value = 11845.5
$ws.Range("K129").Value = 53573202
$ws.Range("L129").Value = 35536.5
$ws.Range("M129").Value = -53568202
$ws.Range("N129").Value = -45536.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2081.25
$ws.Range("I137").Value = 1725.25
$ws.Range("J137").Value = 2437.25
$ws.Range("K137").Value = 5175.75
$ws.Range("L137").Value = 7311.75
$ws.Range("M137").Value = -75.75
$ws.Range("N137").Value = -17511.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 21334.334
$ws.Range("J24").Value = 21334.334
$ws.Range("L24").Value = 21334.334
$ws.Range("N24").Value = -21680.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10699.4
$ws.Range("I70").Value = 10699.4
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 10699.4
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = -10429.4
$ws.Range("M70").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10699.4
$ws.Range("I73").Value = 10699.4
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 10699.4
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = -9763.4
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4923.1665
$ws.Range("I80").Value = 4775.3
$ws.Range("J80").Value = 4997.1
$ws.Range("K80").Value = 4775.3
$ws.Range("L80").Value = 4997.1
$ws.Range("M80").Value = -3777.3
$ws.Range("N80").Value = -6993.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4923.1665
$ws.Range("I83").Value = 4775.3
$ws.Range("J83").Value = 4997.1
$ws.Range("K83").Value = 23876.5
$ws.Range("L83").Value = 24985.5
$ws.Range("M83").Value = -18884.5
$ws.Range("N83").Value = -34969.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3056.027
$ws.Range("I97").Value = 2586.923
$ws.Range("K97").Value = 2586.923
$ws.Range("M97").Value = -2090.923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1044.2
$ws.Range("I107").Value = 1220.2858
$ws.Range("J107").Value = 633.3333
$ws.Range("K107").Value = 1220.2858
$ws.Range("L107").Value = 633.3333
$ws.Range("M107").Value = 699.7141999999999
$ws.Range("N107").Value = -4473.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 14286579
$ws.Range("I16").Value = 17242248
$ws.Range("J16").Value = 846.8333
$ws.Range("K16").Value = 17242248
$ws.Range("L16").Value = 846.8333
$ws.Range("M16").Value = -17242078
$ws.Range("N16").Value = -1186.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 6635.077
$ws.Range("I55").Value = 984.1429000000001
$ws.Range("J55").Value = 13227.833
$ws.Range("K55").Value = 984.1429000000001
$ws.Range("L55").Value = 13227.833
$ws.Range("M55").Value = -811.1429000000001
$ws.Range("N55").Value = -13573.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3250
$ws.Range("I82").Value = 3250
$ws.Range("K82").Value = 3250
$ws.Range("M82").Value = -2889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3250
$ws.Range("I85").Value = 3250
$ws.Range("K85").Value = 3250
$ws.Range("M85").Value = -2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3051.8572
$ws.Range("I132").Value = 1783.7693
$ws.Range("J132").Value = 6715.222
$ws.Range("K132").Value = 5351.3079
$ws.Range("L132").Value = 20145.666
$ws.Range("M132").Value = -2821.3079
$ws.Range("N132").Value = -25205.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19997
$ws.Range("J31").Value = 19997
$ws.Range("L31").Value = 19997
$ws.Range("N31").Value = -20693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5514.5
$ws.Range("I96").Value = 4451
$ws.Range("J96").Value = 6578
$ws.Range("K96").Value = 4451
$ws.Range("L96").Value = 6578
$ws.Range("M96").Value = -3078
$ws.Range("N96").Value = -9324

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 442.5
$ws.Range("I107").Value = 474.83334
$ws.Range("K107").Value = 1424.50002
$ws.Range("M107").Value = 495.4999800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 18666.646
$ws.Range("I122").Value = 2351.9333
$ws.Range("K122").Value = 7055.7999
$ws.Range("M122").Value = -4605.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1377.7188
$ws.Range("I132").Value = 1326.32
$ws.Range("J132").Value = 1561.2858
$ws.Range("K132").Value = 3978.96
$ws.Range("L132").Value = 4683.857400000001
$ws.Range("M132").Value = -1448.96
$ws.Range("N132").Value = -9743.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 371480
$ws.Range("I136").Value = 401127.84
$ws.Range("J136").Value = 882
$ws.Range("K136").Value = 1203383.52
$ws.Range("L136").Value = 2646
$ws.Range("M136").Value = -1200833.52
$ws.Range("N136").Value = -7746

Write-Output "Applied all Lich Profits updates"